$d = $word.ActiveDocument

$replacements = @(
    @("691÷3=230, 1", "549÷9=61, 0"),
    @("136÷9=15, 1", "260÷6=43, 2"),
    @("270÷9=30, 0", "938÷7=134, 0"),
    @("425÷6=70, 5", "678÷7=96, 6"),
    @("988÷8=123, 4", "922÷5=184, 2"),
    @("468÷6=78, 0", "718÷3=239, 1"),
    @("781÷9=86, 7", "986÷5=197, 1"),
    @("437÷5=87, 2", "262÷2=131, 0"),
    @("655÷5=131, 0", "361÷2=180, 1"),
    @("968÷5=193, 3", "552÷3=184, 0"),
    @("275÷5=55, 0", "445÷8=55, 5"),
    @("781÷8=97, 5", "604÷9=67, 1"),
    @("482÷4=120, 2", "441÷6=73, 3"),
    @("407÷9=45, 2", "448÷6=74, 4"),
    @("781÷5=156, 1", "440÷9=48, 8"),
    @("876÷9=97, 3", "263÷7=37, 4"),
    @("718÷8=89, 6", "279÷4=69, 3"),
    @("957÷2=478, 1", "369÷7=52, 5"),
    @("155÷2=77, 1", "120÷7=17, 1"),
    @("793÷3=264, 1", "711÷6=118, 3"),
    @("292÷6=48, 4", "811÷3=270, 1"),
    @("772÷2=386, 0", "362÷2=181, 0"),
    @("362÷7=51, 5", "401÷9=44, 5"),
    @("295÷6=49, 1", "889÷3=296, 1"),
    @("117÷5=23, 2", "711÷6=118, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
